$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D, J, K, L, M, N, O, P, Q for rows 2-4 are cyclically rotated:
# new row2 = old row4, new row3 = old row2, new row4 = old row3
$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

# Capture original values before overwriting
$orig2 = @{}
$orig3 = @{}
$orig4 = @{}
foreach ($col in $cols) {
    $orig2[$col] = $ws.Range("$col" + "2").Value2
    $orig3[$col] = $ws.Range("$col" + "3").Value2
    $orig4[$col] = $ws.Range("$col" + "4").Value2
}

foreach ($col in $cols) {
    $ws.Range("$col" + "2").Value2 = $orig4[$col]
    $ws.Range("$col" + "3").Value2 = $orig2[$col]
    $ws.Range("$col" + "4").Value2 = $orig3[$col]
}
